{"js": "const body = context.document.body;\nconst results = body.search(\"County Court\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'County Court' text in document body.\");\n}\n\n// Append \" Money Claims Centre \" right after \"County Court\" so the\n// paragraph reads \"In the County Court Money Claims Centre \".\nresults.items[0].insertText(\" Money Claims Centre \", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The first table cell reads \"In the County Court\"; rename the court to\n# \"In the County Court Money Claims Centre \" (trailing space retained).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"y Court\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"y Court Money Claims Centre \"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
